# Remove slide 4 ("Inserción de datos con INSERT", sldId 263) together with
# its notes page; PowerPoint automatically renumbers the remaining slides.
$p = $ppt.ActivePresentation
$p.Slides.Item(4).Delete()
